$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = 1.381855808202115
$ws.Range("F1").Value = -1.570796395407675

$ws.Range("E2").Value = 1.382896421261855
$ws.Range("F2").Value = -1.570796393670107

$ws.Range("E3").Value = 1.387559057678701
$ws.Range("F3").Value = -1.570796385884652

$ws.Range("E4").Value = 1.39411894442546
$ws.Range("F4").Value = -1.570796374931256

$ws.Range("E5").Value = 1.398781580842306
$ws.Range("F5").Value = -1.570796367145801

$ws.Range("E6").Value = 1.399822193902047
$ws.Range("F6").Value = -1.570796365408233
